$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'123"
$ws.Range("B2").Value = 0.8709546327590942
$ws.Range("C2").Value = 0.8725749850273132
$ws.Range("D2").Value = 295.7394409179688
$ws.Range("E2").Value = 32.38544464111328
$ws.Range("F2").Value = 33.27213668823242
$ws.Range("G2").Value = 182.0018005371094
$ws.Range("H2").Value = 226.1123199462891
